$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Natmi following Dr Hou advice: a third sending/receiving cluster ("ECs")
# is now included in the LR-pair analysis, expanding the Sending/Target
# cluster cross-product from 2x3 to 3x3 rows and refreshing the derived
# expression/specificity statistics for every row.

# Row 2
$ws.Cells.Item(2, 1).Value = "ECs"
$ws.Cells.Item(2, 2).Value = "Hgf"
$ws.Cells.Item(2, 3).Value = "Met"
$ws.Cells.Item(2, 4).Value = "ECs"
$ws.Cells.Item(2, 5).Value = 2
$ws.Cells.Item(2, 6).Value = 0.6666666666666666
$ws.Cells.Item(2, 7).Value = 5.231719
$ws.Cells.Item(2, 8).Value = 15.695157
$ws.Cells.Item(2, 9).Value = 0.362499186434781
$ws.Cells.Item(2, 10).Value = 0.362499186434781
$ws.Cells.Item(2, 11).Value = 3
$ws.Cells.Item(2, 12).Value = 1
$ws.Cells.Item(2, 13).Value = 1.847798333333333
$ws.Cells.Item(2, 14).Value = 5.543395
$ws.Cells.Item(2, 15).Value = 0.05039680725746681
$ws.Cells.Item(2, 16).Value = 0.05039680725746681
$ws.Cells.Item(2, 17).Value = 9.667161648668333
$ws.Cells.Item(2, 18).Value = 87.00445483801501
$ws.Cells.Item(2, 19).Value = 0.01826880162974218
$ws.Cells.Item(2, 20).Value = 0.01826880162974218

# Row 3
$ws.Cells.Item(3, 1).Value = "ECs"
$ws.Cells.Item(3, 2).Value = "Hgf"
$ws.Cells.Item(3, 3).Value = "Met"
$ws.Cells.Item(3, 4).Value = "FAPs"
$ws.Cells.Item(3, 5).Value = 2
$ws.Cells.Item(3, 6).Value = 0.6666666666666666
$ws.Cells.Item(3, 7).Value = 5.231719
$ws.Cells.Item(3, 8).Value = 15.695157
$ws.Cells.Item(3, 9).Value = 0.362499186434781
$ws.Cells.Item(3, 10).Value = 0.362499186434781
$ws.Cells.Item(3, 11).Value = 3
$ws.Cells.Item(3, 12).Value = 1
$ws.Cells.Item(3, 13).Value = 0.4798556666666667
$ws.Cells.Item(3, 14).Value = 1.439567
$ws.Cells.Item(3, 15).Value = 0.0130875719001099
$ws.Cells.Item(3, 16).Value = 0.0130875719001099
$ws.Cells.Item(3, 17).Value = 2.510470008557667
$ws.Cells.Item(3, 18).Value = 22.594230077019
$ws.Cells.Item(3, 19).Value = 0.004744234166196539
$ws.Cells.Item(3, 20).Value = 0.004744234166196539

# Row 4
$ws.Cells.Item(4, 1).Value = "ECs"
$ws.Cells.Item(4, 2).Value = "Hgf"
$ws.Cells.Item(4, 3).Value = "Met"
$ws.Cells.Item(4, 4).Value = "sCs"
$ws.Cells.Item(4, 5).Value = 2
$ws.Cells.Item(4, 6).Value = 0.6666666666666666
$ws.Cells.Item(4, 7).Value = 5.231719
$ws.Cells.Item(4, 8).Value = 15.695157
$ws.Cells.Item(4, 9).Value = 0.362499186434781
$ws.Cells.Item(4, 10).Value = 0.362499186434781
$ws.Cells.Item(4, 11).Value = 3
$ws.Cells.Item(4, 12).Value = 1
$ws.Cells.Item(4, 13).Value = 34.337334
$ws.Cells.Item(4, 14).Value = 103.012002
$ws.Cells.Item(4, 15).Value = 0.9365156208424232
$ws.Cells.Item(4, 16).Value = 0.9365156208424232
$ws.Cells.Item(4, 17).Value = 179.643282697146
$ws.Cells.Item(4, 18).Value = 1616.789544274314
$ws.Cells.Item(4, 19).Value = 0.3394861506388422
$ws.Cells.Item(4, 20).Value = 0.3394861506388422

# Row 5
$ws.Cells.Item(5, 1).Value = "FAPs"
$ws.Cells.Item(5, 2).Value = "Hgf"
$ws.Cells.Item(5, 3).Value = "Met"
$ws.Cells.Item(5, 4).Value = "ECs"
$ws.Cells.Item(5, 5).Value = 3
$ws.Cells.Item(5, 6).Value = 1
$ws.Cells.Item(5, 7).Value = 8.316945333333335
$ws.Cells.Item(5, 8).Value = 24.950836
$ws.Cells.Item(5, 9).Value = 0.5762706133406404
$ws.Cells.Item(5, 10).Value = 0.5762706133406403
$ws.Cells.Item(5, 11).Value = 3
$ws.Cells.Item(5, 12).Value = 1
$ws.Cells.Item(5, 13).Value = 1.847798333333333
$ws.Cells.Item(5, 14).Value = 5.543395
$ws.Cells.Item(5, 15).Value = 0.05039680725746681
$ws.Cells.Item(5, 16).Value = 0.05039680725746681
$ws.Cells.Item(5, 17).Value = 15.36803772535778
$ws.Cells.Item(5, 18).Value = 138.31233952822
$ws.Cells.Item(5, 19).Value = 0.02904219902867044
$ws.Cells.Item(5, 20).Value = 0.02904219902867043

# Row 6
$ws.Cells.Item(6, 1).Value = "FAPs"
$ws.Cells.Item(6, 2).Value = "Hgf"
$ws.Cells.Item(6, 3).Value = "Met"
$ws.Cells.Item(6, 4).Value = "FAPs"
$ws.Cells.Item(6, 5).Value = 3
$ws.Cells.Item(6, 6).Value = 1
$ws.Cells.Item(6, 7).Value = 8.316945333333335
$ws.Cells.Item(6, 8).Value = 24.950836
$ws.Cells.Item(6, 9).Value = 0.5762706133406404
$ws.Cells.Item(6, 10).Value = 0.5762706133406403
$ws.Cells.Item(6, 11).Value = 3
$ws.Cells.Item(6, 12).Value = 1
$ws.Cells.Item(6, 13).Value = 0.4798556666666667
$ws.Cells.Item(6, 14).Value = 1.439567
$ws.Cells.Item(6, 15).Value = 0.0130875719001099
$ws.Cells.Item(6, 16).Value = 0.0130875719001099
$ws.Cells.Item(6, 17).Value = 3.99093334755689
$ws.Cells.Item(6, 18).Value = 35.918400128012
$ws.Cells.Item(6, 19).Value = 0.007541983086016063
$ws.Cells.Item(6, 20).Value = 0.007541983086016062

# Row 7
$ws.Cells.Item(7, 1).Value = "FAPs"
$ws.Cells.Item(7, 2).Value = "Hgf"
$ws.Cells.Item(7, 3).Value = "Met"
$ws.Cells.Item(7, 4).Value = "sCs"
$ws.Cells.Item(7, 5).Value = 3
$ws.Cells.Item(7, 6).Value = 1
$ws.Cells.Item(7, 7).Value = 8.316945333333335
$ws.Cells.Item(7, 8).Value = 24.950836
$ws.Cells.Item(7, 9).Value = 0.5762706133406404
$ws.Cells.Item(7, 10).Value = 0.5762706133406403
$ws.Cells.Item(7, 11).Value = 3
$ws.Cells.Item(7, 12).Value = 1
$ws.Cells.Item(7, 13).Value = 34.337334
$ws.Cells.Item(7, 14).Value = 103.012002
$ws.Cells.Item(7, 15).Value = 0.9365156208424232
$ws.Cells.Item(7, 16).Value = 0.9365156208424232
$ws.Cells.Item(7, 17).Value = 285.5817297704081
$ws.Cells.Item(7, 18).Value = 2570.235567933672
$ws.Cells.Item(7, 19).Value = 0.5396864312259539
$ws.Cells.Item(7, 20).Value = 0.5396864312259538

# Row 8
$ws.Cells.Item(8, 1).Value = "sCs"
$ws.Cells.Item(8, 2).Value = "Hgf"
$ws.Cells.Item(8, 3).Value = "Met"
$ws.Cells.Item(8, 4).Value = "ECs"
$ws.Cells.Item(8, 5).Value = 3
$ws.Cells.Item(8, 6).Value = 1
$ws.Cells.Item(8, 7).Value = 0.8836963333333333
$ws.Cells.Item(8, 8).Value = 2.651089
$ws.Cells.Item(8, 9).Value = 0.06123020022457864
$ws.Cells.Item(8, 10).Value = 0.06123020022457864
$ws.Cells.Item(8, 11).Value = 3
$ws.Cells.Item(8, 12).Value = 1
$ws.Cells.Item(8, 13).Value = 1.847798333333333
$ws.Cells.Item(8, 14).Value = 5.543395
$ws.Cells.Item(8, 15).Value = 0.05039680725746681
$ws.Cells.Item(8, 16).Value = 0.05039680725746681
$ws.Cells.Item(8, 17).Value = 1.632892611906111
$ws.Cells.Item(8, 18).Value = 14.696033507155
$ws.Cells.Item(8, 19).Value = 0.00308580659905419
$ws.Cells.Item(8, 20).Value = 0.00308580659905419

# Row 9
$ws.Cells.Item(9, 1).Value = "sCs"
$ws.Cells.Item(9, 2).Value = "Hgf"
$ws.Cells.Item(9, 3).Value = "Met"
$ws.Cells.Item(9, 4).Value = "FAPs"
$ws.Cells.Item(9, 5).Value = 3
$ws.Cells.Item(9, 6).Value = 1
$ws.Cells.Item(9, 7).Value = 0.8836963333333333
$ws.Cells.Item(9, 8).Value = 2.651089
$ws.Cells.Item(9, 9).Value = 0.06123020022457864
$ws.Cells.Item(9, 10).Value = 0.06123020022457864
$ws.Cells.Item(9, 11).Value = 3
$ws.Cells.Item(9, 12).Value = 1
$ws.Cells.Item(9, 13).Value = 0.4798556666666667
$ws.Cells.Item(9, 14).Value = 1.439567
$ws.Cells.Item(9, 15).Value = 0.0130875719001099
$ws.Cells.Item(9, 16).Value = 0.0130875719001099
$ws.Cells.Item(9, 17).Value = 0.4240466931625556
$ws.Cells.Item(9, 18).Value = 3.816420238463
$ws.Cells.Item(9, 19).Value = 0.0008013546478972982
$ws.Cells.Item(9, 20).Value = 0.0008013546478972983

# Row 10
$ws.Cells.Item(10, 1).Value = "sCs"
$ws.Cells.Item(10, 2).Value = "Hgf"
$ws.Cells.Item(10, 3).Value = "Met"
$ws.Cells.Item(10, 4).Value = "sCs"
$ws.Cells.Item(10, 5).Value = 3
$ws.Cells.Item(10, 6).Value = 1
$ws.Cells.Item(10, 7).Value = 0.8836963333333333
$ws.Cells.Item(10, 8).Value = 2.651089
$ws.Cells.Item(10, 9).Value = 0.06123020022457864
$ws.Cells.Item(10, 10).Value = 0.06123020022457864
$ws.Cells.Item(10, 11).Value = 3
$ws.Cells.Item(10, 12).Value = 1
$ws.Cells.Item(10, 13).Value = 34.337334
$ws.Cells.Item(10, 14).Value = 103.012002
$ws.Cells.Item(10, 15).Value = 0.9365156208424232
$ws.Cells.Item(10, 16).Value = 0.9365156208424232
$ws.Cells.Item(10, 17).Value = 30.343776152242
$ws.Cells.Item(10, 18).Value = 273.093985370178
$ws.Cells.Item(10, 19).Value = 0.05734303897762714
$ws.Cells.Item(10, 20).Value = 0.05734303897762714

Write-Host "Hgf-Met sheet updated: rows 2-10 refreshed with ECs cluster"
